$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("B1").Value = "is list ordered review"
$ws.Range("C1").Value = "picking right pair review"
$ws.Range("D1").Value = "picking wrong pair review"
$ws.Range("E1").Value = "picking zero pair review"
$ws.Range("F1").Value = "picking only one element review"
$ws.Range("G1").Value = "final score"

# Copy the header style (bold, centered, bordered) from A1 onto the new header cells,
# preserving the cell values already written above (xlPasteFormats = -4122).
$ws.Range("A1").Copy()
$ws.Range("D1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2
$ws.Range("B2").Value = "score:20.0 / 20.0`n"
$ws.Range("C2").Value = "score:4.0 / 20.0`nbruce wayne and wayne enterprises: failed`nclark kent and daily planet: failed`npeter parker and daily bugle: failed`nwillie wonka and chocolate factory: failed"
$ws.Range("D2").Value = "score:0.0 / 20.0`nbruce wayne and chocolate factory: failed"
$ws.Range("E2").Value = "score:20.0 / 20.0`n"
$ws.Range("F2").Value = "score:20.0 / 20.0`n"
$ws.Range("G2").Value = "64.0/100.0"

# Row 3
$ws.Range("B3").Value = "score:20.0 / 20.0`n"
$ws.Range("C3").Value = "score:20.0 / 20.0`n"
$ws.Range("D3").Value = "score:20.0 / 20.0`n"
$ws.Range("E3").Value = "score:20.0 / 20.0`n"
$ws.Range("F3").Value = "score:20.0 / 20.0`n"
$ws.Range("G3").Value = "100.0/100.0"
